$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-10-20 Monday"; New = "2025-10-21 Tuesday" },
    @{ Old = "996÷7="; New = "314÷7=" },
    @{ Old = "794÷8="; New = "454÷9=" },
    @{ Old = "524÷8="; New = "509÷9=" },
    @{ Old = "806÷3="; New = "319÷5=" },
    @{ Old = "330÷9="; New = "919÷3=" },
    @{ Old = "645÷9="; New = "563÷7=" },
    @{ Old = "193÷9="; New = "821÷9=" },
    @{ Old = "540÷9="; New = "271÷4=" },
    @{ Old = "672÷9="; New = "320÷6=" },
    @{ Old = "173÷4="; New = "360÷3=" },
    @{ Old = "958÷5="; New = "456÷2=" },
    @{ Old = "338÷4="; New = "816÷2=" },
    @{ Old = "360÷5="; New = "716÷2=" },
    @{ Old = "200÷2="; New = "709÷6=" },
    @{ Old = "384÷2="; New = "927÷4=" },
    @{ Old = "552÷5="; New = "157÷8=" },
    @{ Old = "122÷3="; New = "770÷6=" },
    @{ Old = "111÷5="; New = "696÷2=" },
    @{ Old = "801÷4="; New = "419÷7=" },
    @{ Old = "404÷5="; New = "474÷6=" },
    @{ Old = "187÷7="; New = "453÷4=" },
    @{ Old = "489÷5="; New = "683÷7=" },
    @{ Old = "461÷5="; New = "497÷4=" },
    @{ Old = "694÷3="; New = "385÷8=" },
    @{ Old = "538÷7="; New = "664÷8=" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)
}
